$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '28.248.52'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +3.50%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.783.74'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.12%  '

# Row 4
$ws.Range('E4').Value = '  +0.17%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '339.88'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.09%  '

# Row 6
$ws.Range('E6').Value = '  +0.17%  '

# Row 7
$ws.Range('E7').Value = '  -3.91%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3436'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.37%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '47.11'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.77%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.156'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -3.17%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07412'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.56%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '23.59'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +8.50%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.003'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.31%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.463'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.02%  '

# Row 15
$ws.Range('B15').NumberFormat = '@'
$ws.Range('B15').Value = 'Chainlink'
$ws.Range('B15').Style = 'Normal'
$ws.Range('C15').NumberFormat = '@'
$ws.Range('C15').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('C15').Style = 'Normal'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.360'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +3.33%  '

# Row 16
$ws.Range('B16').NumberFormat = '@'
$ws.Range('B16').Value = 'WrappedEther'
$ws.Range('B16').Style = 'Normal'
$ws.Range('C16').NumberFormat = '@'
$ws.Range('C16').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('C16').Style = 'Normal'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.799.19'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.72%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001078'
$ws.Range('D17').Style = 'Normal'

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.06698'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.28%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '82.40'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.08%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.002'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.14%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.46'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.30%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.420'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.23%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '28.263.66'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +3.56%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '12.11'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.97%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.365'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.03%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '20.83'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.41%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.421'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -3.53%  '

# Row 28
$ws.Range('B28').NumberFormat = '@'
$ws.Range('B28').Value = 'Monero'
$ws.Range('B28').Style = 'Normal'
$ws.Range('C28').NumberFormat = '@'
$ws.Range('C28').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('C28').Style = 'Normal'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '154.66'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.44%  '

# Row 29
$ws.Range('B29').NumberFormat = '@'
$ws.Range('B29').Value = 'LidoDAOToken'
$ws.Range('B29').Style = 'Normal'
$ws.Range('C29').NumberFormat = '@'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('C29').Style = 'Normal'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.410'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -3.43%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.993.41'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.28%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '135.50'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.17%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.020'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.32%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.117'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.85%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.08959'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.49%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '12.78'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.67%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.02412'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.64%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.6854'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.33%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '5.359'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.94%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.06389'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.67%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.2168'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.52%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.253'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.29%  '

# Row 42
$ws.Range('E42').Value = '  -7.26%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.324'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.17%  '

# Row 44
$ws.Range('B44').NumberFormat = '@'
$ws.Range('B44').Value = 'EnergySwap'
$ws.Range('B44').Style = 'Normal'
$ws.Range('C44').NumberFormat = '@'
$ws.Range('C44').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('C44').Style = 'Normal'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '14.24'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.70%  '

# Row 45
$ws.Range('B45').NumberFormat = '@'
$ws.Range('B45').Value = 'Frax'
$ws.Range('B45').Style = 'Normal'
$ws.Range('C45').NumberFormat = '@'
$ws.Range('C45').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('C45').Style = 'Normal'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.002'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.21%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.6292'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.60%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.884'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.32%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '133.46'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.86%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.080'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -2.57%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.07516'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +5.39%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.196'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +3.65%  '
